# Flip the binary label column (B2:B201) on Sheet2: 0 <-> 1.
# In the original workbook, rows 2-101 (label=0) and rows 102-201 (label=1)
# were the reverse of the desired training-set ordering; this swaps the
# class label for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 2; $i -le 201; $i++) {
    $cell = $ws.Cells.Item($i, 2)
    $current = $cell.Value2
    $cell.Value2 = 1 - $current
}

# Restore the last active selection used when the file was re-saved.
$ws.Range("J8").Select()
